$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) contains text values that often look numeric
# (e.g. "312.36", "45.380.25"). Force the column to Text format before
# writing so Excel does not silently convert these into real numbers,
# then clear the temporary formatting again so the cells end up with
# no explicit style, matching the original workbook's cell styling.
$priceCol = $ws.Range("D2:D51")
$priceCol.NumberFormat = "@"

$ws.Range("D2").Value2 = '45.380.25'
$ws.Range("E2").Value2 = '  -0.36%  '
$ws.Range("D3").Value2 = '2.368.68'
$ws.Range("E3").Value2 = '  -0.78%  '
$ws.Range("D5").Value2 = '312.36'
$ws.Range("E5").Value2 = '  -2.14%  '
$ws.Range("D6").Value2 = '108.15'
$ws.Range("E6").Value2 = '  -3.13%  '
$ws.Range("D7").Value2 = '0.630'
$ws.Range("E7").Value2 = '  -1.27%  '
$ws.Range("E8").Value2 = '  +0.00%  '
$ws.Range("E9").Value2 = '  -3.42%  '
$ws.Range("D10").Value2 = '40.79'
$ws.Range("E10").Value2 = '  -3.11%  '
$ws.Range("E11").Value2 = '  -1.64%  '
$ws.Range("D12").Value2 = '8.45'
$ws.Range("E12").Value2 = '  -2.69%  '
$ws.Range("E13").Value2 = '  +1.18%  '
$ws.Range("E14").Value2 = '  -4.09%  '
$ws.Range("D15").Value2 = '2.729.02'
$ws.Range("E15").Value2 = '  -0.55%  '
$ws.Range("D16").Value2 = '15.31'
$ws.Range("E16").Value2 = '  -2.81%  '
$ws.Range("D17").Value2 = '2.366.73'
$ws.Range("E17").Value2 = '  -1.47%  '
$ws.Range("D18").Value2 = '45.287.49'
$ws.Range("E18").Value2 = '  -0.56%  '
$ws.Range("D19").Value2 = '14.43'
$ws.Range("E19").Value2 = '  +11.22%  '
$ws.Range("B20").Value2 = 'Uniswap'
$ws.Range("C20").Value2 = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").Value2 = '7.24'
$ws.Range("E20").Value2 = '  -5.16%  '
$ws.Range("B21").Value2 = 'ShibaInu'
$ws.Range("C21").Value2 = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D21").Value2 = '0.0000106'
$ws.Range("E21").Value2 = '  -1.90%  '
$ws.Range("E22").Value2 = '  -2.74%  '
$ws.Range("D23").Value2 = '3.55'
$ws.Range("E23").Value2 = '  -0.33%  '
$ws.Range("D24").Value2 = '260.02'
$ws.Range("E24").Value2 = '  -3.45%  '
$ws.Range("E25").Value2 = '  +1.52%  '
$ws.Range("E26").Value2 = '  -0.08%  '
$ws.Range("D27").Value2 = '11.07'
$ws.Range("E27").Value2 = '  -2.40%  '
$ws.Range("D28").Value2 = '7.28'
$ws.Range("E28").Value2 = '  -3.97%  '
$ws.Range("B29").Value2 = 'Hedera'
$ws.Range("C29").Value2 = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D29").Value2 = '0.0973'
$ws.Range("E29").Value2 = '  +2.84%  '
$ws.Range("B30").Value2 = 'Toncoin'
$ws.Range("C30").Value2 = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").Value2 = '2.22'
$ws.Range("E30").Value2 = '  -4.64%  '
$ws.Range("E31").Value2 = '  -3.05%  '
$ws.Range("D32").Value2 = '36.80'
$ws.Range("E32").Value2 = '  -4.86%  '
$ws.Range("D33").Value2 = '166.62'
$ws.Range("E33").Value2 = '  -1.95%  '
$ws.Range("E34").Value2 = '  -3.86%  '
$ws.Range("E35").Value2 = '  -2.39%  '
$ws.Range("E36").Value2 = '  -1.15%  '
$ws.Range("D37").Value2 = '4.69'
$ws.Range("E37").Value2 = '  -3.78%  '
$ws.Range("E38").Value2 = '  +8.56%  '
$ws.Range("D39").Value2 = '3.94'
$ws.Range("E39").Value2 = '  +0.26%  '
$ws.Range("E40").Value2 = '  -4.51%  '
$ws.Range("E41").Value2 = '  -3.51%  '
$ws.Range("D42").Value2 = '98.50'
$ws.Range("E42").Value2 = '  -7.27%  '
$ws.Range("D43").Value2 = '69.93'
$ws.Range("E43").Value2 = '  -2.17%  '
$ws.Range("D44").Value2 = '0.227'
$ws.Range("E44").Value2 = '  -5.76%  '
$ws.Range("D45").Value2 = '12.82'
$ws.Range("E45").Value2 = '  -7.55%  '
$ws.Range("D46").Value2 = '0.999'
$ws.Range("E46").Value2 = '  -0.23%  '
$ws.Range("D47").Value2 = '1.819.51'
$ws.Range("E47").Value2 = '  +9.98%  '
$ws.Range("D48").Value2 = '5.86'
$ws.Range("E48").Value2 = '  +0.83%  '
$ws.Range("D49").Value2 = '83.23'
$ws.Range("E49").Value2 = '  +5.22%  '
$ws.Range("D50").Value2 = '110.85'
$ws.Range("E50").Value2 = '  -6.31%  '
$ws.Range("E51").Value2 = '  -0.89%  '

$priceCol.ClearFormats()
